# Updated symbol list on Fri Dec 16 06:43:03 UTC 2022 with GitHub Actions
# Refreshes the crypto price/volume table on Sheet1: new price quotes for
# existing rows, plus the coin-ranking list shifting up one slot (row 10's
# "One" entry rolls off and a new "One" entry reappears at row 18).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''261.32'

# Row 3
$ws.Range('D3').Value = '''22.88'

# Row 4
$ws.Range('D4').Value = '''6.183'

# Row 5
$ws.Range('D5').Value = '''0.06240'

# Row 6
$ws.Range('D6').Value = '''6.724'

# Row 8
$ws.Range('D8').Value = '''1.341'

# Row 9
$ws.Range('D9').Value = '''0.7972'

# Row 10
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '''0.1576'
$ws.Range('E10').Value = '9WazirXWRX'

# Row 11
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '''0.08113'
$ws.Range('E11').Value = '10MandalaExchangeTokenMDX'

# Row 12
$ws.Range('B12').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C12').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D12').Value = '''0.03474'
$ws.Range('E12').Value = '11LiechtensteinCryptoassetsExchangeLCX'

# Row 13
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D13').Value = '''0.03089'
$ws.Range('E13').Value = '12BitrueCoinBTR'

# Row 14
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D14').Value = '''0.09320'
$ws.Range('E14').Value = '13BitMartTokenBMX'

# Row 15
$ws.Range('B15').Value = 'MCDex'
$ws.Range('C15').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('D15').Value = '''3.845'
$ws.Range('E15').Value = '14MCDexMCB'

# Row 16
$ws.Range('B16').Value = 'BitForexToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D16').Value = '''0.001712'
$ws.Range('E16').Value = '15BitForexTokenBF'

# Row 17
$ws.Range('B17').Value = 'CoinExToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('D17').Value = '''0.04779'
$ws.Range('E17').Value = '16CoinExTokenCET'

# Row 18
$ws.Range('B18').Value = 'One'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D18').Value = '''0.0006140'
$ws.Range('E18').Value = '17OneONEWorstin24h'

# Row 19
$ws.Range('D19').Value = '''0.006203'

# Row 20
$ws.Range('D20').Value = '''0.001091'

# Row 21
$ws.Range('D21').Value = '''0.004081'

# Row 23
$ws.Range('D23').Value = '''3.696'

# Row 24
$ws.Range('D24').Value = '''2.212'

# Row 25
$ws.Range('D25').Value = '''0.3363'

# Row 27
$ws.Range('D27').Value = '''0.0003202'

# Row 40
$ws.Range('D40').Value = '''0.04607'

# Row 41
$ws.Range('D41').Value = '''0.007139'

# Row 42
$ws.Range('D42').Value = '''0.1118'

# Row 43
$ws.Range('D43').Value = '''0.003130'

# Row 45
$ws.Range('D45').Value = '''0.002970'

# Row 46
$ws.Range('D46').Value = '''0.00005889'

# Row 48
$ws.Range('D48').Value = '''0.7000'

# Row 49
$ws.Range('D49').Value = '''0.1098'
$ws.Range('E49').Value = '48BOLOBOLO'

# Row 50
$ws.Range('D50').Value = '''0.00002100'
